$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TODO note inserted as D38 (set first so it becomes the new shared
# string placed immediately before the "Indie Project..." string, matching
# the author's edit order where the new si precedes the edited one).
$ws.Range("D38").Value = "TODO - consider different approach for recording story removals -- I have a bad hunch about current approach"

# Week 7 entry: append "and super basic entity classes" to the existing note,
# and bump the logged hours for that day from 3.5 to 4.
$ws.Range("D36").Value = "Indie Project: organized priorities, figured out how to generate UML diagrams from IntelliJ, tried and failed to recreate Log4J problem, set up Servlet shells and super basic entity classes`nWeek 7: Intro video"
$ws.Range("B36").Value = 4

# Update the view so the newly-added row is in frame and selected, matching
# the scrolled/selected state left behind by the edit.
[void]$ws.Range("D58").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 36
$win.ScrollColumn = 1
